# Adds new category items to the bottom of "sheet".
# Mirrors the authored diff: 5 new rows (Others, Cold Pressed Juice,
# Pooja Items, Stationary, Printing Material) appended after the
# existing data, with the previous "last block" of rows re-based onto
# the plain border-only style so the new rows become the visually
# distinct trailing block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 314
$newItems = @(
    "Others",
    "Cold Pressed Juice",
    "Pooja Items",
    "Stationary",
    "Printing Material"
)

# Re-format the previous tail block (rows 309-314) to match the plain
# border style used earlier in the column (same look, used before the
# new trailing block was introduced).
$formatDonor = $ws.Range("A304")
$formatDonor.Copy()
$ws.Range("A309:A314").PasteSpecial(-4122)

# Append the new rows, each carrying on the distinct trailing style
# (copied from what used to be the last formatted cell).
$styleDonor = $ws.Range("A309")
for ($i = 0; $i -lt $newItems.Count; $i++) {
    $targetRow = $lastRow + 1 + $i
    $cell = $ws.Cells.Item($targetRow, 1)
    $styleDonor.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $newItems[$i]
}

$ws.Application.CutCopyMode = $false

# Leave the active selection on the final appended cell, matching the
# saved workbook's cursor position.
$ws.Range("A" + ($lastRow + $newItems.Count)).Select()
